$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JUNE-22")

# Row 11 and 10's new unique text is entered first (D11 before D10) so that the
# shared-strings table ends up built in the same append order as the source
# workbook: 208=Rlogic support note, 209=sony task note, 210=PR-Details note,
# 211=IE->Edge migration note.
$ws.Range("C11").Value = "RPA Rlogic"
$ws.Range("D11").Value = "2. Supported to Rlogic Extended warranty task as records counts are missmatching, now an additional logs are recorded in the log file and the task is rerunning"

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = 44720
$ws.Range("B10").NumberFormat = "m/d/yy"
$ws.Range("C10").Value = "RPA SONY"
$ws.Range("D10").Value = "1. For the sony task,  doing some Research and development for clicking Ok button to select certificate to enter into sony site and it is fixed by Mohan san, whereas the task testing is pending.`n"
$ws.Range("E10").Value = 0.4
$ws.Range("E10").NumberFormat = "0%"
$ws.Range("F10").Value = "WIP"

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = 44721
$ws.Range("B12").NumberFormat = "m/d/yy"
$ws.Range("C12").Value = "RPA GSS"
$ws.Range("D12").Value = "1. Correction received for the PR-Details task due to task failed during upload, and it has been fixed, tested and it is running smoothly from download to upload."
$ws.Range("E12").Value = 1
$ws.Range("E12").NumberFormat = "0%"
$ws.Range("F12").Value = "Completed"

$ws.Range("C13").Value = "RPA SONY"
$ws.Range("D13").Value = "2. Migration from IE to Edge browser task  has been completed for the SONY_APNEWSIS_DAILY_V2 task, tested and it is running smoothly from download to upload `n(Inbound daily task, outbound task, and RPSI inquiry task includes)"
$ws.Range("E13").Value = 1
$ws.Range("E13").NumberFormat = "0%"
$ws.Range("F13").Value = "Completed"

$ws.Range("D10:D13").WrapText = $true

$ws.Rows("10").RowHeight = 43.2
$ws.Rows("11").RowHeight = 28.8
$ws.Rows("12").RowHeight = 28.8
$ws.Rows("13").RowHeight = 43.2

$ws.Range("D13").Select()
